$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "proxy4" worksheet after the last existing sheet (proxy3).
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "proxy4"

# ---------------------------------------------------------------------------
# 2. Header row (re-uses existing shared strings: age, d11B, d11Bsd, d18O, MgCa, species)
# ---------------------------------------------------------------------------
$newSheet.Range("A1").Value = "age"
$newSheet.Range("B1").Value = "d11B"
$newSheet.Range("C1").Value = "d11Bsd"
$newSheet.Range("D1").Value = "d18O"
$newSheet.Range("E1").Value = "MgCa"
$newSheet.Range("F1").Value = "species"

# ---------------------------------------------------------------------------
# 3. Data rows 2-14
# ---------------------------------------------------------------------------
$newSheet.Range("A2").Value = 58.776000000000003
$newSheet.Range("B2").Value = 16.670000000000002
$newSheet.Range("C2").Value = 0.13
$newSheet.Range("D2").Value = -1.1299999999999999
$newSheet.Range("E2").Value = 3.59
$newSheet.Range("F2").Value = "Grub"

$newSheet.Range("A3").Value = 58.621000000000002
$newSheet.Range("D3").Value = -1.1299999999999999
$newSheet.Range("E3").Value = 3.5
$newSheet.Range("F3").Value = "Grub"

$newSheet.Range("A4").Value = 58.101999999999997
$newSheet.Range("D4").Value = -1.59
$newSheet.Range("F4").Value = "Grub"

$newSheet.Range("A5").Value = 58.100999999999999
$newSheet.Range("B5").Value = 16.670000000000002
$newSheet.Range("C5").Value = 0.115
$newSheet.Range("E5").Value = 4.4000000000000004
$newSheet.Range("F5").Value = "Grub"

$newSheet.Range("A6").Value = 58.085999999999999
$newSheet.Range("E6").Value = 4.47
$newSheet.Range("F6").Value = "Grub"

$newSheet.Range("A7").Value = 57.987000000000002
$newSheet.Range("B7").Value = 16.47
$newSheet.Range("C7").Value = 0.13
$newSheet.Range("D7").Value = -1.65
$newSheet.Range("E7").Value = 3.88
$newSheet.Range("F7").Value = "Grub"

$newSheet.Range("A8").Value = 57.225999999999999
$newSheet.Range("E8").Value = 3.58
$newSheet.Range("F8").Value = "Grub"

$newSheet.Range("A9").Value = 57.073999999999998
$newSheet.Range("D9").Value = -1.42
$newSheet.Range("E9").Value = 3.77
$newSheet.Range("F9").Value = "Grub"

$newSheet.Range("A10").Value = 55.95
$newSheet.Range("B10").Value = 15.73
$newSheet.Range("C10").Value = 0.19500000000000001
$newSheet.Range("D10").Value = -1.53
$newSheet.Range("E10").Value = 3.33
$newSheet.Range("F10").Value = "Grub"

$newSheet.Range("A11").Value = 55.932000000000002
$newSheet.Range("D11").Value = -1.99
$newSheet.Range("E11").Value = 5.04
$newSheet.Range("F11").Value = "Grub"

$newSheet.Range("A12").Value = 55.884999999999998
$newSheet.Range("B12").Value = 14.9
$newSheet.Range("C12").Value = 0.15
$newSheet.Range("D12").Value = -1.84
$newSheet.Range("F12").Value = "Grub"

$newSheet.Range("A13").Value = 55.787999999999997
$newSheet.Range("D13").Value = -1.57
$newSheet.Range("E13").Value = 3.98
$newSheet.Range("F13").Value = "Grub"

$newSheet.Range("A14").Value = 55.738999999999997
$newSheet.Range("E14").Value = 3.84
$newSheet.Range("F14").Value = "Grub"

# ---------------------------------------------------------------------------
# 4. Update proxy3's view: select A1:F41 (was A19:XFD19) and leave it as the
#    non-active sheet.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("proxy3")
$ws3.Activate()
$ws3.Range("A1:F41").Select()

# ---------------------------------------------------------------------------
# 5. Make proxy4 the active sheet/tab with the selection on D24.
# ---------------------------------------------------------------------------
$newSheet.Activate()
$newSheet.Range("D24").Select()
